# Generate Report for Handoff
# Adds a new tracked file (83c2c437-a725-460d-8aab-5787596344fc.md) as row 3
# on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the existing
# row 2 entry for 7f2803df-b3bd-43ca-b3ca-1fe2ea6e28c1.md.

$wb = $excel.ActiveWorkbook

$guid = "83c2c437-a725-460d-8aab-5787596344fc"
$hash = "65470f45d57efe315b15fde9624c7e04dda10ea4"
$repoUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c83fdd8b6a371053f89995b0b970c3a7355c9c7e/e2e/$guid.md"
$displayPath = "e2e\$guid.md"

# ---------------------------------------------------------------
# Overview sheet (row 3)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = "$guid.md"
$ws.Range("B3").Value = $displayPath
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-22 16:42:22"
$ws.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Add($ws.Range("B3"), $repoUrl, $null, $null, $displayPath)
$ws.Range("B3").Font.Color = 15570276
$ws.Range("B3").Font.Underline = $true

# ---------------------------------------------------------------
# zh-cn sheet (row 3)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3").Value = "$guid.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "$guid.$hash.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-22 16:42:18"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "'"
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$ws.Hyperlinks.Add($ws.Range("A3"), $repoUrl, $null, $null, "$guid.md")
$ws.Range("A3").Font.Color = 15570276
$ws.Range("A3").Font.Underline = $true

# ---------------------------------------------------------------
# de-de sheet (row 3)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3").Value = "$guid.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "$guid.$hash.de-de.xlf"
$ws.Range("H3").Value = "2016-08-22 16:42:22"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "'"
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

$ws.Hyperlinks.Add($ws.Range("A3"), $repoUrl, $null, $null, "$guid.md")
$ws.Range("A3").Font.Color = 15570276
$ws.Range("A3").Font.Underline = $true

Write-Host "Added handoff row for $guid"
